# Scheduled runner: refresh market-board derived profit figures on the
# Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Re-pricing changed
# currentAveragePrice(NQ/HQ) inputs, which cascade into LevePrice/LeveProfit
# columns for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1719.6
$ws.Range("J97").Value = 1719.6
$ws.Range("L97").Value = 5158.799999999999
$ws.Range("N97").Value = -6150.799999999999
$ws.Range("H132").Value = 1111
$ws.Range("I132").Value = 1111
$ws.Range("K132").Value = 3333
$ws.Range("M132").Value = -803
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 10226.3
$ws.Range("J138").Value = 8140.3335
$ws.Range("L138").Value = 24421.0005
$ws.Range("N138").Value = -34701.00049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1402.3334
$ws.Range("I61").Value = 1373.5
$ws.Range("J61").Value = 1460
$ws.Range("K61").Value = 1373.5
$ws.Range("L61").Value = 1460
$ws.Range("M61").Value = -1161.5
$ws.Range("N61").Value = -1884
$ws.Range("H86").Value = 25000
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27372
$ws.Range("H89").Value = 25000
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86856
$ws.Range("H136").Value = 1402.3334
$ws.Range("I136").Value = 1373.5
$ws.Range("J136").Value = 1460
$ws.Range("K136").Value = 4120.5
$ws.Range("L136").Value = 4380
$ws.Range("M136").Value = -1570.5
$ws.Range("N136").Value = -9480
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 1771
$ws.Range("I80").Value = 1066.6666
$ws.Range("J80").Value = 2299.25
$ws.Range("K80").Value = 1066.6666
$ws.Range("L80").Value = 2299.25
$ws.Range("M80").Value = -68.66660000000002
$ws.Range("N80").Value = -4295.25
$ws.Range("H83").Value = 1771
$ws.Range("I83").Value = 1066.6666
$ws.Range("J83").Value = 2299.25
$ws.Range("K83").Value = 5333.333000000001
$ws.Range("L83").Value = 11496.25
$ws.Range("M83").Value = -341.3330000000005
$ws.Range("N83").Value = -21480.25
$ws.Range("H86").Value = 4860.6313
$ws.Range("I86").Value = 5060.857
$ws.Range("J86").Value = 4300
$ws.Range("K86").Value = 5060.857
$ws.Range("L86").Value = 4300
$ws.Range("M86").Value = -3937.857
$ws.Range("N86").Value = -6546
$ws.Range("H89").Value = 4860.6313
$ws.Range("I89").Value = 5060.857
$ws.Range("J89").Value = 4300
$ws.Range("K89").Value = 25304.285
$ws.Range("L89").Value = 21500
$ws.Range("M89").Value = -19688.285
$ws.Range("N89").Value = -32732
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3097.1667
$ws.Range("I94").Value = 2438.5
$ws.Range("K94").Value = 2438.5
$ws.Range("M94").Value = -1987.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 18000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3946
$ws.Range("I132").Value = 3946
$ws.Range("K132").Value = 11838
$ws.Range("M132").Value = -9308
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5148.3335
$ws.Range("I16").Value = 4222.5
$ws.Range("K16").Value = 4222.5
$ws.Range("M16").Value = -4052.5
$ws.Range("H22").Value = 2139.1
$ws.Range("I22").Value = 1986.375
$ws.Range("K22").Value = 1986.375
$ws.Range("M22").Value = -1691.375
$ws.Range("H27").Value = 2139.1
$ws.Range("I27").Value = 1986.375
$ws.Range("K27").Value = 1986.375
$ws.Range("M27").Value = -1879.375
$ws.Range("H46").Value = 4666.6665
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H88").Value = 69420
$ws.Range("J88").Value = 69420
$ws.Range("L88").Value = 69420
$ws.Range("N88").Value = -70276
$ws.Range("H91").Value = 69420
$ws.Range("J91").Value = 69420
$ws.Range("L91").Value = 69420
$ws.Range("N91").Value = -72384
$ws.Range("H93").Value = 1239.4
$ws.Range("I93").Value = 1174.25
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1174.25
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 73.75
$ws.Range("N93").Value = -3996
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 4500
$ws.Range("N132").Value = -9560
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3445
$ws.Range("I81").Value = 3445
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6890
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5829
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 3445
$ws.Range("I84").Value = 3445
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 34450
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -29146
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 1076.2142
$ws.Range("I107").Value = 995.7143
$ws.Range("K107").Value = 2987.1429
$ws.Range("M107").Value = -1067.1429
$ws.Range("H132").Value = 2249.8
$ws.Range("I132").Value = 2166.6667
$ws.Range("K132").Value = 6500.000100000001
$ws.Range("M132").Value = -3970.000100000001
